$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("A6").Value = "9fe83a0f-e574-4337-9aa4-fa11e0a8075d"
$ws.Range("B6").Value = -1.293631102816713
$ws.Range("C6").Value = 36.80742457755191
# Force text, not date, interpretation of the date-like string, then
# restore the default "Normal" style so no extra formatting is left behind.
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2025-10-28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "10:26:28 AM"
$ws.Range("F6").Value = "::1"
$ws.Range("G6").Value = "test"

# Row 7
$ws.Range("A7").Value = "01735368-36c0-43b2-a7f8-3af44c2a1dff"
$ws.Range("B7").Value = -1.2936585497144812
$ws.Range("C7").Value = 36.80746346087199
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "2025-10-28"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "10:34:00 AM"
$ws.Range("F7").Value = "::1"
$ws.Range("G7").Value = "brian"

# Extend the "number stored as text" ignored-error suppression to cover
# the newly added rows (A1:G5 -> A1:G7), matching the rest of the table.
$ws.Range("A1:G7").Errors.Item(9).Ignore = $true
